# Update levequest profit-tracking figures per the scheduled-runner refresh.
# Each sheet is a Table_<ClassAbbr> (Leve Name/Item/.../LeveProfitHQ); only the
# price/profit columns (H-N) for specific rows change - no structural/header edits.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 25003124
$ws.Range("I74").Value = 25003124
$ws.Range("K74").Value = 25003124
$ws.Range("M74").Value = -25002188
# Row 77
$ws.Range("H77").Value = 25003124
$ws.Range("I77").Value = 25003124
$ws.Range("K77").Value = 125015620
$ws.Range("M77").Value = -125010940
# Row 116
$ws.Range("H116").Value = 560981.9399999999
$ws.Range("I116").Value = 1251684.9
$ws.Range("K116").Value = 1251684.9
$ws.Range("M116").Value = -1248242.9

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10353.517
$ws.Range("I32").Value = 6848.2065
$ws.Range("J32").Value = 17714.666
$ws.Range("K32").Value = 6848.2065
$ws.Range("L32").Value = 17714.666
$ws.Range("M32").Value = -6561.2065
$ws.Range("N32").Value = -18288.666
# Row 35
$ws.Range("H35").Value = 21666.25
$ws.Range("I35").Value = 10000
$ws.Range("J35").Value = 33332.5
$ws.Range("K35").Value = 10000
$ws.Range("L35").Value = 33332.5
$ws.Range("M35").Value = -9594
$ws.Range("N35").Value = -34144.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 269.44
$ws.Range("I80").Value = 124.4
$ws.Range("J80").Value = 366.13333
$ws.Range("K80").Value = 124.4
$ws.Range("L80").Value = 366.13333
$ws.Range("M80").Value = 873.6
$ws.Range("N80").Value = -2362.13333
# Row 83
$ws.Range("H83").Value = 269.44
$ws.Range("I83").Value = 124.4
$ws.Range("J83").Value = 366.13333
$ws.Range("K83").Value = 622
$ws.Range("L83").Value = 1830.66665
$ws.Range("M83").Value = 4370
$ws.Range("N83").Value = -11814.66665
# Row 105
$ws.Range("H105").Value = 33335634
$ws.Range("I105").Value = 47621050
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 47621050
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -47619303
$ws.Range("N105").Value = -6494

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4338.385
$ws.Range("I31").Value = 1872.6428
$ws.Range("J31").Value = 7215.0835
$ws.Range("K31").Value = 1872.6428
$ws.Range("L31").Value = 7215.0835
$ws.Range("M31").Value = -1577.6428
$ws.Range("N31").Value = -7805.0835
# Row 34
$ws.Range("H34").Value = 4338.385
$ws.Range("I34").Value = 1872.6428
$ws.Range("J34").Value = 7215.0835
$ws.Range("K34").Value = 1872.6428
$ws.Range("L34").Value = 7215.0835
$ws.Range("M34").Value = -1670.6428
$ws.Range("N34").Value = -7619.0835
# Row 55
$ws.Range("H55").Value = 15000
$ws.Range("I55").Value = 15000
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 15000
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -14685
$ws.Range("N55").ClearContents()
# Row 60
$ws.Range("H60").Value = 27029.666
$ws.Range("J60").Value = 28388.928
$ws.Range("L60").Value = 28388.928
$ws.Range("N60").Value = -29410.928

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 129
$ws.Range("H129").Value = 2419.775
$ws.Range("I129").Value = 2383.125
$ws.Range("K129").Value = 7149.375
$ws.Range("M129").Value = -2149.375
# Row 133
$ws.Range("H133").Value = 3397.7273
$ws.Range("J133").Value = 3040
$ws.Range("L133").Value = 9120
$ws.Range("N133").Value = -19240
# Row 134
$ws.Range("H134").Value = 3390.7827
$ws.Range("I134").Value = 2434.8572
$ws.Range("J134").Value = 4877.778
$ws.Range("K134").Value = 7304.571599999999
$ws.Range("L134").Value = 14633.334
$ws.Range("M134").Value = -2234.571599999999
$ws.Range("N134").Value = -24773.334
# Row 137
$ws.Range("H137").Value = 9686.357
$ws.Range("I137").Value = 2739.2307
$ws.Range("J137").Value = 99999
$ws.Range("K137").Value = 8217.6921
$ws.Range("L137").Value = 299997
$ws.Range("M137").Value = -3117.6921
$ws.Range("N137").Value = -310197
# Row 138
$ws.Range("H138").Value = 2991.3845
$ws.Range("I138").Value = 2061.25
$ws.Range("J138").Value = 4479.6
$ws.Range("K138").Value = 6183.75
$ws.Range("L138").Value = 13438.8
$ws.Range("M138").Value = -1043.75
$ws.Range("N138").Value = -23718.8
# Row 139
$ws.Range("H139").Value = 3152
$ws.Range("I139").Value = 3152
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 9456
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -4316
$ws.Range("N139").ClearContents()
# Row 141
$ws.Range("H141").Value = 8789.909
$ws.Range("I141").Value = 8336.125
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 25008.375
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -19828.375
$ws.Range("N141").Value = -40360

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 52
$ws.Range("H52").Value = 36000
$ws.Range("J52").Value = 36000
$ws.Range("L52").Value = 36000
$ws.Range("N52").Value = -36518
# Row 53
$ws.Range("H53").Value = 29997
$ws.Range("J53").Value = 29997
$ws.Range("L53").Value = 29997
$ws.Range("N53").Value = -31259
# Row 70
$ws.Range("H70").Value = 5057.6523
$ws.Range("I70").Value = 4785.56
$ws.Range("J70").Value = 5381.5713
$ws.Range("K70").Value = 4785.56
$ws.Range("L70").Value = 5381.5713
$ws.Range("M70").Value = -4515.56
$ws.Range("N70").Value = -5921.5713
# Row 73
$ws.Range("H73").Value = 5057.6523
$ws.Range("I73").Value = 4785.56
$ws.Range("J73").Value = 5381.5713
$ws.Range("K73").Value = 4785.56
$ws.Range("L73").Value = 5381.5713
$ws.Range("M73").Value = -3849.56
$ws.Range("N73").Value = -7253.5713

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 20300
$ws.Range("J64").Value = 20300
$ws.Range("L64").Value = 20300
$ws.Range("N64").Value = -20796
# Row 67
$ws.Range("H67").Value = 20300
$ws.Range("J67").Value = 20300
$ws.Range("L67").Value = 20300
$ws.Range("N67").Value = -22016
# Row 110
$ws.Range("H110").Value = 39750
$ws.Range("J110").Value = 39750
$ws.Range("L110").Value = 39750
$ws.Range("N110").Value = -47930
# Row 111
$ws.Range("H111").Value = 39800
$ws.Range("J111").Value = 39800
$ws.Range("L111").Value = 39800
$ws.Range("N111").Value = -47980
# Row 112
$ws.Range("H112").Value = 39750
$ws.Range("J112").Value = 39750
$ws.Range("L112").Value = 39750
$ws.Range("N112").Value = -42704
# Row 113
$ws.Range("H113").Value = 10375.2
$ws.Range("I113").Value = 20300.8
$ws.Range("J113").Value = 449.6
$ws.Range("K113").Value = 60902.39999999999
$ws.Range("L113").Value = 1348.8
$ws.Range("M113").Value = -58732.39999999999
$ws.Range("N113").Value = -5688.8
# Row 125
$ws.Range("H125").Value = 41612.855
$ws.Range("J125").Value = 41612.855
$ws.Range("L125").Value = 41612.855
$ws.Range("N125").Value = -51452.855
